# Fixed naive component forecaster bug - Presentation state 11.02.
# Recomputed the YoY forecast vectors with the corrected naive component
# forecaster: C2 (a stray duplicate y_1 value) is removed, and every
# y_1 / y_1_forecast value is refreshed with its corrected (float-noise
# level) result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 was an errant leftover cell - clear it out entirely so the
# <c> element is removed, then refresh E2 with the recomputed value.
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 7.865470614547321

# Row 3
$ws.Range("E3").Value = -6.760862998203621

# Row 4
$ws.Range("C4").Value = 0.5799958470386724

# Row 6
$ws.Range("C6").Value = 0.5930547804883446
$ws.Range("E6").Value = -1.194610791899986

# Row 8
$ws.Range("E8").Value = 7.617133650412167

# Row 9
$ws.Range("C9").Value = 1.670328650030162
$ws.Range("E9").Value = 2.037906845818593

# Row 10
$ws.Range("C10").Value = 2.562791874943349

# Row 11
$ws.Range("C11").Value = 1.526411006965578
$ws.Range("E11").Value = 0.6601843988560452

# Row 12
$ws.Range("E12").Value = 1.55185774637272

# Row 14
$ws.Range("E14").Value = -5.866344937500012

# Row 15
$ws.Range("C15").Value = -2.616267413525608
$ws.Range("E15").Value = -4.982381489483368

# Row 17
$ws.Range("C17").Value = -1.298607950737285

# Row 18
$ws.Range("C18").Value = -0.994151974263302

# Row 19
$ws.Range("C19").Value = 1.069485063776932
$ws.Range("E19").Value = -2.110726282892139
